# sejours.xlsx edit script
# Corrections + ajout de nouveaux tarifs en parallele des anciens
#
# Summary of changes:
#  - Split "proprietaire" column into "Prenom propriétaire" (new) + "nom proprietaire" (renamed)
#  - Insert a new row 2 with long-form human-readable descriptions for every column
#    (yellow highlight, word-wrap, bold/plain Calibri 11)
#  - Resize several columns
#  - Misc cosmetic: window position, selection, page setup (portrait, paper size 9)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new column before column I to host "Prenom propriétaire".
#    The previous column I ("proprietaire") shifts to column J and is
#    renamed to "nom proprietaire" (now holding only the last name).
# ---------------------------------------------------------------------------
$ws.Columns.Item(9).EntireColumn.Insert()

# Give the new header cell (I1) the same look as the rest of row 1 (bold, filled)
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("I1").Value2 = "Prénom propriétaire"
$ws.Range("J1").Value2 = "nom proprietaire"

# ---------------------------------------------------------------------------
# 2. Insert a new descriptive row right under the header row. The former
#    data row (row 2) shifts down to row 3.
# ---------------------------------------------------------------------------
$ws.Rows.Item(2).EntireRow.Insert()

$ws.Range("A2").Value2 = "Champ technique, ne pas remplir"
$ws.Range("B2").Value2 = "Date d'arrivée en pension pour ce séjour, au format JJ/MM/AAAA HH:MM"
$ws.Range("C2").Value2 = "Date de départ de pension pour ce séjour, au format JJ/MM/AAAA HH:MM"
$ws.Range("D2").Value2 = "Nombre de cages fournies par le propriétaire"
$ws.Range("E2").Value2 = "Nombre de cage à fournir par la pension"
$ws.Range("F2").Value2 = "Montant total du séjour en pension"
$ws.Range("G2").Value2 = "Montant restant à payer pour ce séjour"
$ws.Range("H2").Value2 = "Animaux du séjour avec exactement le même nom que dans l'import des animaux, séparés par une virgule s'il y en a plusieurs. Si possible mettre en évidence les animaux qui pourraient poser problème à cause de doublons"
$ws.Range("I2").Value2 = "Prénom du propriétaire"
$ws.Range("J2").Value2 = "Nom du propriétaire"
$ws.Range("K2").Value2 = "L'animal devra-t-il/aura-t-il été vacciné par la pension pendant son séjour? OUI ou NON"
$ws.Range("L2").Value2 = "Un des animaux necessite t'il des soins particuliers pendant le séjour? OUI ou NON"
$ws.Range("M2").Value2 = "Les soins particuliers de l'animal se font-ils par injection? OUI ou NON"
$ws.Range("N2").Value2 = "Commentaire éventuel sur le séjour."

$descRange = $ws.Range("A2:N2")
$descRange.Font.Bold = $false
$descRange.Interior.Color = 10092543   # RGB(255,255,153) -> FFFFFF99, yellow highlight
$descRange.WrapText = $true
$ws.Rows.Item(2).RowHeight = 101.5

# ---------------------------------------------------------------------------
# 3. Column widths
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 18.084    # A  -> 19
$ws.Columns.Item(2).ColumnWidth = 27.336    # B  -> 28.1796875
$ws.Columns.Item(3).ColumnWidth = 26.418    # C  -> 27.26953125
$ws.Columns.Item(4).ColumnWidth = 17.084    # D  -> 18
$ws.Columns.Item(5).ColumnWidth = 25.25     # E  -> 26.08984375
$ws.Columns.Item(6).ColumnWidth = 18.25     # F  -> 19.08984375
$ws.Columns.Item(8).ColumnWidth = 31.418    # H  -> 32.26953125
$ws.Columns.Item(9).ColumnWidth = 18.918    # I  -> 19.81640625
$ws.Columns.Item(10).ColumnWidth = 16.584   # J  -> 17.453125
$ws.Columns.Item(11).ColumnWidth = 17.584   # K  -> 18.453125
$ws.Columns.Item(12).ColumnWidth = 15.418   # L  -> 16.36328125
$ws.Columns.Item(13).ColumnWidth = 14.584   # M  -> 15.54296875

# ---------------------------------------------------------------------------
# 4. Cosmetics: window position, active selection, page setup
# ---------------------------------------------------------------------------
$excel.ActiveWindow.Left = -110

[void]$ws.Range("H2").Select()

$ws.PageSetup.PaperSize = 9      # xlPaperA4
$ws.PageSetup.Orientation = 1    # xlPortrait
